$d = $word.ActiveDocument

# --- Paragraph 1: the **ID__...__ID** marker paragraph ---
$p1 = $d.Paragraphs(1)

# 1) Add a paragraph border (top/left/bottom/right) with 5pt "space" on each side,
#    matching <w:pBdr><w:top w:space="5"/>...</w:pBdr>.
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# 2) Change the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# 3) Remove the trailing run that contains only a literal space
#    (runs: [0,31) = "**ID__AFFARS_5352_topic_2__ID**", [31,32) = " ").
$trailingSpace = $d.Range(31, 32)
$trailingSpace.Delete()

# 4) Update the marker text itself (now the only run left in the paragraph).
$null = $p1.Range.Find.Execute("**ID__AFFARS_5352_topic_2__ID**", $true, $false, $false,
                                $false, $false, $true, 1, $false,
                                "**ID__AFFARS_SUBPART_5352_2__ID**", 2)
